$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update a handful of per-country statistics (COVID data refresh) ---

$lastRow = 212

for ($r = 4; $r -le $lastRow; $r++) {
    $name = $ws.Cells.Item($r, 1).Text

    if ($name -eq "Reino Unido") {
        $ws.Cells.Item($r, 6).Value = 1559
    }
    elseif ($name -eq "Israel") {
        $ws.Cells.Item($r, 2).Value = 8018
        $ws.Cells.Item($r, 3).Value = 167
        $ws.Cells.Item($r, 5).Value = 7545
        $ws.Cells.Item($r, 6).Value = 127
        $ws.Cells.Item($r, 7).Value = 2
        $ws.Cells.Item($r, 8).Value = 46
    }
    elseif ($name -eq "Pakistan") {
        $ws.Cells.Item($r, 2).Value = 2880
        $ws.Cells.Item($r, 3).Value = 62
        $ws.Cells.Item($r, 4).Value = 170
        $ws.Cells.Item($r, 5).Value = 2665
        $ws.Cells.Item($r, 6).Value = 18
        $ws.Cells.Item($r, 7).Value = 4
        $ws.Cells.Item($r, 8).Value = 45
    }
    elseif ($name -eq "Hungria") {
        $ws.Cells.Item($r, 2).Value = 733
        $ws.Cells.Item($r, 3).Value = 55
        $ws.Cells.Item($r, 4).Value = 66
        $ws.Cells.Item($r, 5).Value = 633
        $ws.Cells.Item($r, 6).Value = 17
        $ws.Cells.Item($r, 7).Value = 2
        $ws.Cells.Item($r, 8).Value = 34
    }
    elseif ($name -eq "Bulgaria") {
        $ws.Cells.Item($r, 2).Value = 522
        $ws.Cells.Item($r, 3).Value = 19
        $ws.Cells.Item($r, 4).Value = 37
        $ws.Cells.Item($r, 5).Value = 467
        $ws.Cells.Item($r, 6).Value = 26
        $ws.Cells.Item($r, 7).Value = 1
        $ws.Cells.Item($r, 8).Value = 18
    }
}

# Updated totals can change country ranking - keep the table sorted by
# "Casos totales" (column B) descending, as it was before the edit.
$sortRange = $ws.Range("A4:H" + $lastRow)
$sortKey = $ws.Range("B4")
$sortRange.Sort($sortKey, 2)

# --- Update the "last updated" footer text ---
$ws.Range("A1").Value = "Datos actualizados a 5 de Abril de 2020 a las 07:52"
